# Fruta / hortaliza, semanal
#
# A new weekly price-report row is inserted at row 322 (pushing the
# existing rows 322-340 down to 323-341). The new row carries a fresh
# observation (Región de Arica y Parinacota, 2023-08-09) while every
# other row keeps its original data, just shifted down by one position.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 322..340 down to 323..341, leaving row 322 empty.
$ws.Rows.Item(322).Insert()

# Populate the newly inserted row 322 with the new observation.
$ws.Range("A322").Value = 7
$ws.Range("B322").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C322").Value = "Ñuble"
$ws.Range("D322").Value = 45147
$ws.Range("E322").Value = 16
$ws.Range("F322").Value = 100112032
$ws.Range("G322").Value = "Zapallo italiano"
$ws.Range("H322").Value = "Sin especificar"
$ws.Range("I322").Value = "Primera"
$ws.Range("J322").Value = 120
$ws.Range("K322").Value = 14000
$ws.Range("L322").Value = 15000
$ws.Range("M322").Value = 14667
$ws.Range("N322").Value = "`$/caja 50 unidades"
$ws.Range("O322").Value = "Región de Arica y Parinacota"
$ws.Range("P322").Value = 293
$ws.Range("Q322").Value = 50
$ws.Range("R322").Value = "Hortaliza"
